$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cookie")

# --- Insert a new column, landing it at F, while keeping the exact custom width ---
# Inserting directly before F would lose float precision on the custom column width
# (the COM ColumnWidth setter quantizes to 1/7-character steps), so instead we insert
# inside the existing D:E custom-width block at E (which inherits the exact width
# from its neighbours with no precision loss), then move the old column E's contents
# back into E and leave the brand new, still-identically-widthed column at F.
$ws.Columns("E").Insert()
$ws.Range("F1:F10").Copy($ws.Range("E1:E10"))
$ws.Range("F1:F10").ClearContents()

# --- New "InitSoulStone" column header/type row ---
$ws.Range("F1").Value2 = "InitSoulStone"
$ws.Range("F2").Value2 = "int"

# --- New InitSoulStone data values (rows 3-10) ---
$ws.Range("F3").Value2 = 20
$ws.Range("F4").Value2 = 20
$ws.Range("F5").Value2 = 20
$ws.Range("F6").Value2 = 20
$ws.Range("F7").Value2 = 20
$ws.Range("F8").Value2 = 20
$ws.Range("F9").Value2 = 20
$ws.Range("F10").Value2 = 20

# --- Updated values in column G (previously column F before the insert) ---
$ws.Range("G8").Value2 = 20
$ws.Range("G9").Value2 = 20

# --- Updated "Num" (column A) values ---
$ws.Range("A3").Value2 = 1010
$ws.Range("A4").Value2 = 1020
$ws.Range("A5").Value2 = 2010
$ws.Range("A6").Value2 = 3010
$ws.Range("A7").Value2 = 4010
$ws.Range("A8").Value2 = 5010
$ws.Range("A9").Value2 = 6010
$ws.Range("A10").Value2 = 6020

# --- Restore the selection shown in the workbook ---
$ws.Range("F2").Select()
